$wb = $excel.ActiveWorkbook

# --- About sheet updates ---
$about = $wb.Worksheets.Item("About")

# Row 10 (A10): update "Note" continuation text
$about.Range("A10").Value2 = "GWP values for VOC, CO, and NOx vary by region.  These values are"

# Row 11 (A11): update text to the new India-specific note
$about.Range("A11").Value2 = 'updated for India from the IPCC source, using "South Asia" rows in the tables.'

# Row 12 is no longer needed (its text merged/removed) - delete the entire row,
# which shifts row 14 ("BC and OC values only include...") up to row 13.
$about.Rows.Item(12).Delete()

# --- Data sheet updates ---
$data = $wb.Worksheets.Item("Data")

# VOC row (row 3): North America Row 3 -> South Asia Row 4, values updated
$data.Range("B3").Value2 = 27.8
$data.Range("C3").Value2 = 8.8
$data.Range("D3").Value2 = "p. 740, Table 8.A.5, Row 4 (VOC South Asia)"

# CO row (row 4): North America Row 3 -> South Asia Row 4, values updated
$data.Range("B4").Value2 = 5.7
$data.Range("C4").Value2 = 1.8
$data.Range("D4").Value2 = "p. 740, Table 8.A.4, Row 4 (CO South Asia)"

# NOx row (row 5): North America Row 3 -> South Asia Row 4, values updated
$data.Range("B5").Value2 = -40.7
$data.Range("C5").Value2 = -25.3
$data.Range("D5").Value2 = "p. 739, Table 8.A.3, Row 4 (NOx South Asia)"

$wb.Save()
